# Project-plan sheet update:
#  - "个人主页" (personal-home) rows get split out into more granular design /
#    front-end tasks, and new back-end / database rows are added
#    ("include font-awesome and complete personal-home").
#  - Rows 11-12 (old trailing rows) are removed, the table now ends at row 10.
#  - B9 loses its border (style goes from the bordered cell style to the
#    plain centered one) while getting its new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- rewrite the data rows (row 1 header is unchanged) ---------------------
$ws.Range("B2").Value = "具体图片"
$ws.Range("C2").Value = "设计"

$ws.Range("B3").Value = "广场页面"
$ws.Range("C3").Value = "设计"

$ws.Range("B4").Value = "个人设置"
$ws.Range("C4").Value = "设计"

$ws.Range("B5").Value = "个人设置"
$ws.Range("C5").Value = "前台实现"
$ws.Range("D5").Value = 0

$ws.Range("B6").Value = "广场页面"
$ws.Range("C6").Value = "前台实现"

$ws.Range("B7").Value = "具体图片"
$ws.Range("C7").Value = "前台实现"
$ws.Range("D7").Value = 1

$ws.Range("B8").Value = "总体"
$ws.Range("C8").Value = "MVC构架"

$ws.Range("B9").Value = "登录注册"
$ws.Range("C9").Value = "后台实现"
$ws.Range("D9").Value = 2
$ws.Range("B9").Borders.LineStyle = -4142   # xlLineStyleNone - drop the border on B9

$ws.Range("B10").Value = "数据库"
$ws.Range("C10").Value = "数据库构架"
$ws.Range("D10").Value = 2

# --- drop the old rows 11 and 12 (table now ends at row 10) ----------------
$ws.Rows("11:12").Delete()

# --- move the active selection like the author left it ---------------------
$null = $ws.Range("F9").Select()
